# Trade #82 closed at 2026-02-17 15:52:59 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.96
$wsSummary.Range("B4").Value = -0.05
$wsSummary.Range("B6").Value = 82
$wsSummary.Range("B8").Value = 43
$wsSummary.Range("B9").Value = 32.93

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.95999999999999
$wsStatus.Range("D4").Value = 82
$wsStatus.Range("E4").Value = -0.05
$wsStatus.Range("F4").Value = -0.04
$wsStatus.Range("G4").Value = 32.93

# ---- New trade row (#82) appended to "All Trades" and "MarketMaking" sheets ----
$newRow = @{
    A = 82
    B = "2026-02-17"
    C = "15:52:52"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.189263
    G = 0.18
    H = "CLOSED"
    I = -4.894
    J = -0.01
    K = 99.95999999999999
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A83").Value = $newRow.A
    # Force the date-looking string to stay text (matches the existing
    # "Date" column cells above, which are plain text, not real dates).
    $ws.Range("B83").NumberFormat = "@"
    $ws.Range("B83").Value = $newRow.B
    $ws.Range("C83").Value = $newRow.C
    $ws.Range("D83").Value = $newRow.D
    $ws.Range("E83").Value = $newRow.E
    $ws.Range("F83").Value = $newRow.F
    $ws.Range("G83").Value = $newRow.G
    $ws.Range("H83").Value = $newRow.H
    $ws.Range("I83").Value = $newRow.I
    $ws.Range("J83").Value = $newRow.J
    $ws.Range("K83").Value = $newRow.K
    $ws.Range("L83").Value = $newRow.L
    $ws.Range("M83").Value = $newRow.M
    $ws.Range("N83").Value = $newRow.N
    $ws.Range("O83").Value = $newRow.O
    $ws.Range("P83").Value = $newRow.P
    $ws.Range("Q83").Value = $newRow.Q
}
